$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MW CNN 1D")
# ("MW CNN 1D" is also the ActiveSheet / 3rd tab in this workbook.)

# New column group headers for the "210 Packets" run (U:W)
$ws.Range("V1").Value = "210 Packets"
$ws.Range("U2").Value = "Acc"
$ws.Range("V2").Value = "Loss"
$ws.Range("W2").Value = "Time"

# Fill in the 230 Packets run data (Q:S), rows 3-51
$ws.Range("Q3").Value = 87.703657150268498
$ws.Range("R3").Value = 0.32339541383366199
$ws.Range("S3").Value = 513.13632917404095
$ws.Range("Q4").Value = 89.245146512985201
$ws.Range("R4").Value = 0.254406929885792
$ws.Range("S4").Value = 502.66289806365899
$ws.Range("Q5").Value = 89.223802089691105
$ws.Range("R5").Value = 0.259136820126531
$ws.Range("S5").Value = 501.67518973350502
$ws.Range("Q6").Value = 88.185071945190401
$ws.Range("R6").Value = 0.27531552605310899
$ws.Range("S6").Value = 498.93015694618202
$ws.Range("Q7").Value = 89.026963710784898
$ws.Range("R7").Value = 0.249530878661437
$ws.Range("S7").Value = 499.93049550056401
$ws.Range("Q8").Value = 89.589017629623399
$ws.Range("R8").Value = 0.273849304244485
$ws.Range("S8").Value = 502.08425998687699
$ws.Range("Q9").Value = 89.854627847671495
$ws.Range("R9").Value = 0.23938763019855799
$ws.Range("S9").Value = 501.63039469718899
$ws.Range("Q10").Value = 87.841200828552203
$ws.Range("R10").Value = 0.28207105664682902
$ws.Range("S10").Value = 499.93942928314198
$ws.Range("Q11").Value = 88.208788633346501
$ws.Range("R11").Value = 0.29210190150156001
$ws.Range("S11").Value = 504.15216732025101
$ws.Range("Q12").Value = 88.374793529510498
$ws.Range("R12").Value = 0.25190459219327799
$ws.Range("S12").Value = 502.88548517227099
$ws.Range("Q13").Value = 89.603245258331299
$ws.Range("R13").Value = 0.280288274066054
$ws.Range("S13").Value = 503.613765001297
$ws.Range("Q14").Value = 89.152652025222693
$ws.Range("R14").Value = 0.25122237992315799
$ws.Range("S14").Value = 503.34356117248501
$ws.Range("Q15").Value = 88.955819606781006
$ws.Range("R15").Value = 0.262904541021033
$ws.Range("S15").Value = 501.38673090934702
$ws.Range("Q16").Value = 89.155024290084796
$ws.Range("R16").Value = 0.25923464412610497
$ws.Range("S16").Value = 503.03025245666498
$ws.Range("Q17").Value = 88.346338272094698
$ws.Range("R17").Value = 0.26661843286817999
$ws.Range("S17").Value = 502.25094246864302
$ws.Range("Q18").Value = 89.173996448516803
$ws.Range("R18").Value = 0.25083581785206399
$ws.Range("S18").Value = 503.68373632431002
$ws.Range("Q19").Value = 88.789814710616994
$ws.Range("R19").Value = 0.263333906400956
$ws.Range("S19").Value = 502.88818001747097
$ws.Range("Q20").Value = 87.5803351402282
$ws.Range("R20").Value = 0.32525435395212399
$ws.Range("S20").Value = 499.70321536064102
$ws.Range("Q21").Value = 88.853842020034705
$ws.Range("R21").Value = 0.28285851866833001
$ws.Range("S21").Value = 503.92856574058499
$ws.Range("Q22").Value = 88.353449106216402
$ws.Range("R22").Value = 0.28909618343189503
$ws.Range("S22").Value = 505.04830193519501
$ws.Range("Q23").Value = 89.173996448516803
$ws.Range("R23").Value = 0.269333901143752
$ws.Range("S23").Value = 503.45758676528902
$ws.Range("Q24").Value = 89.399296045303302
$ws.Range("R24").Value = 0.26802358747653798
$ws.Range("S24").Value = 503.64614963531398
$ws.Range("Q25").Value = 87.777173519134493
$ws.Range("R25").Value = 0.31135509708749098
$ws.Range("S25").Value = 504.09301447868302
$ws.Range("Q26").Value = 88.237249851226807
$ws.Range("R26").Value = 0.294473698137734
$ws.Range("S26").Value = 503.687465190887
$ws.Range("Q27").Value = 88.882303237914996
$ws.Range("R27").Value = 0.26118712204175798
$ws.Range("S27").Value = 505.29672431945801
$ws.Range("Q28").Value = 85.941612720489502
$ws.Range("R28").Value = 0.31737491708827198
$ws.Range("S28").Value = 504.03100633621199
$ws.Range("Q29").Value = 87.881517410278306
$ws.Range("R29").Value = 0.30841567893434002
$ws.Range("S29").Value = 504.01990795135498
$ws.Range("Q30").Value = 87.385869026183997
$ws.Range("R30").Value = 0.32520304055804999
$ws.Range("S30").Value = 503.82376790046601
$ws.Range("Q31").Value = 88.545542955398503
$ws.Range("R31").Value = 0.28942292421129601
$ws.Range("S31").Value = 501.82972598075798
$ws.Range("Q32").Value = 89.057791233062702
$ws.Range("R32").Value = 0.24997637207201101
$ws.Range("S32").Value = 502.504124641418
$ws.Range("Q33").Value = 88.201671838760305
$ws.Range("R33").Value = 0.29460028703839403
$ws.Range("S33").Value = 501.97546529769897
$ws.Range("Q34").Value = 88.813525438308702
$ws.Range("R34").Value = 0.251193278323396
$ws.Range("S34").Value = 502.832363128662
$ws.Range("Q35").Value = 88.970047235488806
$ws.Range("R35").Value = 0.26526913727101797
$ws.Range("S35").Value = 500.71220993995598
$ws.Range("Q36").Value = 89.259374141693101
$ws.Range("R36").Value = 0.27159849401903202
$ws.Range("S36").Value = 502.38794183731
$ws.Range("Q37").Value = 89.667272567748995
$ws.Range("R37").Value = 0.25650285274602402
$ws.Range("S37").Value = 501.87644219398499
$ws.Range("Q38").Value = 88.818269968032794
$ws.Range("R38").Value = 0.26652290204589202
$ws.Range("S38").Value = 501.35289955139098
$ws.Range("Q39").Value = 89.157396554946899
$ws.Range("R39").Value = 0.244437376629909
$ws.Range("S39").Value = 502.79593753814697
$ws.Range("Q40").Value = 89.036452770233097
$ws.Range("R40").Value = 0.254452926719276
$ws.Range("S40").Value = 504.56968164443902
$ws.Range("Q41").Value = 88.327366113662706
$ws.Range("R41").Value = 0.28876673847038897
$ws.Range("S41").Value = 502.76340556144697
$ws.Range("Q42").Value = 87.639623880386296
$ws.Range("R42").Value = 0.30783324543025598
$ws.Range("S42").Value = 507.10781574249199
$ws.Range("Q43").Value = 88.858586549758897
$ws.Range("R43").Value = 0.274864964552544
$ws.Range("S43").Value = 503.01429700851401
$ws.Range("Q44").Value = 88.467282056808401
$ws.Range("R44").Value = 0.27505873285733401
$ws.Range("S44").Value = 505.76502108573902
$ws.Range("Q45").Value = 88.830125331878605
$ws.Range("R45").Value = 0.25426941915770801
$ws.Range("S45").Value = 506.609098911285
$ws.Range("Q46").Value = 89.211940765380803
$ws.Range("R46").Value = 0.26151681250084002
$ws.Range("S46").Value = 504.03819203376702
$ws.Range("Q47").Value = 88.732898235321002
$ws.Range("R47").Value = 0.28200692942259398
$ws.Range("S47").Value = 501.86926364898602
$ws.Range("Q48").Value = 89.209568500518799
$ws.Range("R48").Value = 0.25169849239608999
$ws.Range("S48").Value = 504.42156839370699
$ws.Range("Q49").Value = 88.977164030075002
$ws.Range("R49").Value = 0.267249805656951
$ws.Range("S49").Value = 503.68649840354902
$ws.Range("Q50").Value = 87.786656618118201
$ws.Range("R50").Value = 0.28814849022857197
$ws.Range("S50").Value = 504.78978967666598
$ws.Range("Q51").Value = 88.616693019866901
$ws.Range("R51").Value = 0.26249624012110501
$ws.Range("S51").Value = 505.93578052520701

$ws.Range("L16").Select()
